$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking / percent-looking values are written as literal text
# (matching the workbook convention where these columns store text, not numbers)
$textCells = @(
    "D2",
    "E2",
    "D3",
    "E3",
    "D4",
    "E4",
    "D5",
    "E5",
    "E6",
    "D7",
    "E7",
    "D8",
    "E8",
    "D9",
    "E9",
    "D10",
    "E10",
    "D11",
    "E11",
    "D12",
    "E12",
    "D13",
    "E13",
    "D14",
    "E14",
    "D15",
    "E15",
    "D16",
    "E16",
    "D17",
    "E17",
    "D18",
    "E18",
    "E19",
    "D20",
    "E20",
    "D21",
    "E21",
    "D22",
    "E22",
    "D23",
    "E23",
    "D24",
    "E24",
    "D25",
    "E25",
    "E26",
    "D38",
    "E38",
    "D39",
    "E39",
    "D40",
    "E40",
    "E41",
    "D42",
    "E42",
    "D43",
    "E43",
    "D44",
    "E44",
    "D45",
    "D46",
    "E46",
    "D47",
    "E47",
    "D48",
    "E48",
    "D49",
    "E49",
    "D50",
    "E50"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = "329.25"
$ws.Range("E2").Value = "0.64%"
$ws.Range("D3").Value = "44.11"
$ws.Range("E3").Value = "0.83%"
$ws.Range("D4").Value = "5.579"
$ws.Range("E4").Value = "2.11%"
$ws.Range("D5").Value = "0.08081"
$ws.Range("E5").Value = "0.12%"
$ws.Range("E6").Value = "5.12%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "4.332"
$ws.Range("E7").Value = "0.88%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "0.9522"
$ws.Range("E8").Value = "1.64%"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").Value = "2.569"
$ws.Range("E9").Value = "-4.87%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "0.1162"
$ws.Range("E10").Value = "-2.39%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "0.1856"
$ws.Range("E11").Value = "-2.19%"
$ws.Range("B12").Value = "MCDex"
$ws.Range("C12").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D12").Value = "11.87"
$ws.Range("E12").Value = "37.68%"
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").Value = "0.09846"
$ws.Range("E13").Value = "3.82%"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "0.04744"
$ws.Range("E14").Value = "14.17%"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "0.1068"
$ws.Range("E15").Value = "0.07%"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "0.001286"
$ws.Range("E16").Value = "0.48%"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "0.04234"
$ws.Range("E17").Value = "-2.95%"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").Value = "0.005988"
$ws.Range("E18").Value = "-1.05%"
$ws.Range("E19").Value = "-5.65%"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "0.3472"
$ws.Range("E20").Value = "-0.75%"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "0.1410"
$ws.Range("E21").Value = "3.17%"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").Value = "0.2508"
$ws.Range("E22").Value = "0.48%"
$ws.Range("B23").Value = "BitKan"
$ws.Range("C23").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D23").Value = "0.001256"
$ws.Range("E23").Value = "1.86%"
$ws.Range("B24").Value = "HotbitToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D24").Value = "0.004352"
$ws.Range("E24").Value = "0.23%"
$ws.Range("D25").Value = "0.0001190"
$ws.Range("E25").Value = "-3.50%"
$ws.Range("E26").Value = "-0.55%"
$ws.Range("D38").Value = "0.02658"
$ws.Range("E38").Value = "-0.02%"
$ws.Range("D39").Value = "0.05551"
$ws.Range("E39").Value = "2.48%"
$ws.Range("D40").Value = "0.007565"
$ws.Range("E40").Value = "-0.89%"
$ws.Range("E41").Value = "1.34%"
$ws.Range("D42").Value = "0.008086"
$ws.Range("E42").Value = "-22.59%"
$ws.Range("D43").Value = "0.002016"
$ws.Range("E43").Value = "-3.93%"
$ws.Range("D44").Value = "0.008896"
$ws.Range("E44").Value = "-8.33%"
$ws.Range("D45").Value = "0.00007255"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("E46").Value = "-0.26%"
$ws.Range("B47").Value = "CoinbaseStockToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D47").Value = "0.002272"
$ws.Range("E47").Value = "-0.19%"
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").Value = "0.004726"
$ws.Range("E48").Value = "32.73%"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("E49").Value = "-0.26%"
$ws.Range("D50").Value = "0.0002001"
$ws.Range("E50").Value = "-0.26%"
